$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sensor readings (columns C:H) shift up by one logical sample: the reading
# previously recorded for the next timestamp now belongs to the current one,
# i.e. the original first sample (timestamp 0) is dropped and the series is
# extended with eleven newly captured samples at the end (timestamps
# 1900-2900). Columns A (timestamp) and B (label) stay a simple fixed
# sequence, so every row of the table (A2:H31) is rewritten explicitly here
# to land on the exact target values.
$data = @(
    @(0, "walkingToRunning", -3.284394264221191, -6.207508087158203, -9.360366821289062, -0.09578663945773691, 0.5251570686625937, -0.1342098090959656),
    @(100, "walkingToRunning", -2.543862819671631, -2.628618955612183, -11.02495956420898, -1.14914536605709, 1.190383828204603, 1.292009428791357),
    @(200, "walkingToRunning", 9.523374557495115, -6.183863162994385, -2.546638011932373, -0.09556020285196953, 2.255030672907248, 2.840824183058616),
    @(300, "walkingToRunning", 10.87078380584717, -0.2485504150390625, -1.897834300994873, -1.178172554082913, 0.5742749858305702, 1.787655900641919),
    @(400, "walkingToRunning", -25.65192222595215, -16.10989189147949, -27.1768741607666, -0.017849953278261, 0.6819284758058108, 2.102289744045419),
    @(500, "walkingToRunning", -17.75581932067871, -3.263657808303833, -17.43449020385742, -0.2933185870808658, 0.3665561543570639, 1.667365985216152),
    @(600, "walkingToRunning", 21.72449111938477, -16.76241874694824, -7.28642463684082, -0.2884348296888843, 2.223420211586359, 1.218851031312623),
    @(700, "walkingToRunning", 3.684173583984375, -16.01864433288574, -9.982423782348633, 0.8913533052672417, -1.240040208982344, -4.165942465049586),
    @(800, "walkingToRunning", -0.0730928182601928, 6.555614948272705, -3.539213180541992, 0.04892210568783151, -2.769314762474837, -5.229776274178921),
    @(900, "walkingToRunning", 1.041769504547119, -30.36503982543945, 31.26595878601075, 1.79423694823672, 0.2591471291970247, -2.070258497615964),
    @(1000, "walkingToRunning", 38.83012771606445, 12.28195858001709, -19.02286911010743, 0.327913294667793, 3.156352582185162, 1.962593454381687),
    @(1100, "walkingToRunning", -7.558793067932129, -12.19493198394775, -16.61202049255371, -5.980557255122894, 6.358928622830835, 1.785391607146397),
    @(1200, "walkingToRunning", 0.2677345275878906, -73.31625366210938, 24.66749572753906, -5.223221169577771, 5.128131230672261, -2.261452886793327),
    @(1300, "walkingToRunning", 1.331151604652405, 2.387452363967896, -8.010488510131836, 2.6634052724078, -1.262619167134376, -3.230092231778134),
    @(1400, "walkingToRunning", 16.03682136535645, -33.019287109375, 43.97341537475586, 4.784483177074534, -2.167020964737715, -4.068033393454396),
    @(1500, "walkingToRunning", 28.34896659851075, 6.394547939300537, -27.16730880737305, -0.6180673411502886, -0.3114872364606223, -3.968297384787306),
    @(1600, "walkingToRunning", 2.31614351272583, -4.862353324890137, -18.97574234008789, -6.191343093263846, -1.469820827677687, 4.923526774281952),
    @(1700, "walkingToRunning", -14.28880977630615, -23.93070220947266, 2.869308471679688, -6.191955239300739, 0.06456651607004327, 5.562430750344849),
    @(1800, "walkingToRunning", 6.518649578094482, -0.0271482467651367, 0.3105573654174804, 7.060775899656733, -5.641036204093896, -6.224707568901178),
    @(1900, "walkingToRunning", 9.343793869018556, 9.33461856842041, 8.671773910522461, 2.895918472953405, -3.924228730409128, 1.001562802687943),
    @(2000, "walkingToRunning", 16.0465145111084, -45.62339782714844, 13.98112106323242, 6.419859485349798, 3.571867339277037, 0.7840970496912497),
    @(2100, "walkingToRunning", 18.59218597412109, 6.586655139923096, -14.55167388916016, -3.988556474879168, 3.858240864703046, 3.390691111053246),
    @(2200, "walkingToRunning", -13.51971435546875, -18.57223892211914, -24.33260536193848, -2.792354390241082, 2.638935013093731, 1.279087337894736),
    @(2300, "walkingToRunning", -12.54771614074707, -18.61154365539551, 28.26399230957031, -4.14442459281516, 9.648339557186993, -0.673063791604435),
    @(2400, "walkingToRunning", -1.716351509094239, 5.821096897125244, -9.188434600830078, -1.790030563511108, 9.842405033572287, -8.527358793406183),
    @(2500, "walkingToRunning", 41.91740417480469, -71.66004180908203, 11.51219272613525, 6.688283160112809, -5.523453207983309, -0.2542450458244243),
    @(2600, "walkingToRunning", 22.04729652404785, 0.830233097076416, -27.47162818908692, 3.419822825325905, -5.20332591199646, -2.518984390917574),
    @(2700, "walkingToRunning", -15.56076526641846, -8.609291076660156, -21.82845687866211, -0.0826715539042695, 0.4171818759705612, 2.113873891784312),
    @(2800, "walkingToRunning", -10.87422180175781, -26.33984756469727, 9.715606689453123, -12.35876123801507, -14.20579128680002, 4.713338022646673),
    @(2900, "walkingToRunning", 4.591959953308105, 5.822259902954102, -12.29100227355957, 5.30543631401633, -16.49609409092704, 8.442643778335675)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
